{"js": "// Add a centered, bold title block (\"IDS 720 \u2013 Team 8 \u2013 Opioid Policy\n// Effects\" / \"Project Outline\") plus a blank list-style spacer paragraph\n// to the very top of the document body, matching the author's\n// \"Minor changes before submit\" edit.\n\nconst body = context.document.body;\n\n// Insert in reverse order at the start so the final order reads:\n// 1) Title line 1\n// 2) Title line 2\n// 3) Blank spacer paragraph\n// 4) ... existing content ...\n\nconst spacer = body.insertParagraph(\"\", Word.InsertLocation.start);\nspacer.style = \"List Paragraph\";\nspacer.leftIndent = 18; // 360 twips = 18 pt\n\nconst titleLine2 = body.insertParagraph(\"Project Outline\", Word.InsertLocation.start);\ntitleLine2.style = \"Normal\";\ntitleLine2.alignment = Word.Alignment.centered;\ntitleLine2.leftIndent = 18; // 360 twips = 18 pt\ntitleLine2.firstLineIndent = -18; // -360 twips => w:hanging=\"360\"\ntitleLine2.font.bold = true;\ntitleLine2.font.boldBidirectional = true;\n\nconst titleLine1 = body.insertParagraph(\"IDS 720 \\u2013 Team 8 \\u2013 Opioid Policy Effects\", Word.InsertLocation.start);\ntitleLine1.style = \"Normal\";\ntitleLine1.alignment = Word.Alignment.centered;\ntitleLine1.leftIndent = 18; // 360 twips = 18 pt\ntitleLine1.firstLineIndent = -18; // -360 twips => w:hanging=\"360\"\ntitleLine1.font.bold = true;\ntitleLine1.font.boldBidirectional = true;\n\nawait context.sync();\n", "ps1": "# Add a centered, bold title block (\"IDS 720 \u2013 Team 8 \u2013 Opioid Policy\n# Effects\" / \"Project Outline\") plus a blank list-style spacer paragraph\n# to the very top of the document body, matching the author's\n# \"Minor changes before submit\" edit.\n\n$d = $word.ActiveDocument\n\n# Insert the two title lines as a single operation at the very start of\n# the document (this will initially inherit the first paragraph's\n# numbered-list formatting, which we strip off below).\n$r = $d.Range(0, 0)\n$r.InsertBefore(\"IDS 720 \u2013 Team 8 \u2013 Opioid Policy Effects`rProject Outline`r\")\n\n$titleLine1 = $d.Paragraphs(1)\n$titleLine1.Range.ListFormat.RemoveNumbers()\n$titleLine1.Style = \"Normal\"\n$titleLine1.Range.Font.Bold = 1\n$titleLine1.Range.Font.BoldBi = 1\n$titleLine1.Alignment = 1\n$titleLine1.LeftIndent = 18\n$titleLine1.FirstLineIndent = -18\n\n$titleLine2 = $d.Paragraphs(2)\n$titleLine2.Range.ListFormat.RemoveNumbers()\n$titleLine2.Style = \"Normal\"\n$titleLine2.Range.Font.Bold = 1\n$titleLine2.Range.Font.BoldBi = 1\n$titleLine2.Alignment = 1\n$titleLine2.LeftIndent = 18\n$titleLine2.FirstLineIndent = -18\n\n# Blank spacer paragraph (List Paragraph style, left-indented) between the\n# title block and the existing \"Topic:\" line.\n$spacerRange = $d.Paragraphs(3).Range\n$spacerRange.Collapse(1)\n$spacerRange.InsertBefore(\"`r\")\n$spacer = $d.Paragraphs(3)\n$spacer.Range.ListFormat.RemoveNumbers()\n$spacer.Style = \"List Paragraph\"\n$spacer.LeftIndent = 18\n"}
